# Reorder the "Recorded By" (column G) entries so that "System" is moved
# to the end of the comma-separated list of recorders, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#   "system, System, backup@backdoor.com" -> "system, backup@backdoor.com, System"
# Rows where "System" is not present, or is already the last entry, are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ",\s*"
    $hasSystem = $false
    $newParts = @()
    foreach ($part in $parts) {
        if ($part.Equals("System")) {
            $hasSystem = $true
        } else {
            $newParts += $part
        }
    }

    if ($hasSystem) {
        $newParts += "System"
        $newValue = $newParts -join ", "
        if ($newValue -ne $value) {
            $cell.Value = $newValue
        }
    }
}
